$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1139.75
$ws.Range("I9").Value = 1404
$ws.Range("J9").Value = 875.5
$ws.Range("K9").Value = 1404
$ws.Range("L9").Value = 875.5
$ws.Range("M9").Value = -1235
$ws.Range("N9").Value = -1213.5
$ws.Range("H31").Value = 1799.75
$ws.Range("I31").Value = 1508.8182
$ws.Range("K31").Value = 4526.4546
$ws.Range("M31").Value = -4296.4546
$ws.Range("H33").Value = 2482.8462
$ws.Range("I33").Value = 2850.5
$ws.Range("J33").Value = 460.75
$ws.Range("K33").Value = 2850.5
$ws.Range("L33").Value = 460.75
$ws.Range("M33").Value = -2621.5
$ws.Range("N33").Value = -918.75
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H62").Value = 8969.5
$ws.Range("I62").Value = 8969.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 8969.5
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -8345.5
$ws.Range("H65").Value = 8969.5
$ws.Range("I65").Value = 8969.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 44847.5
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -41727.5
$ws.Range("H74").Value = 10159.267
$ws.Range("I74").Value = 3125
$ws.Range("J74").Value = 14848.777
$ws.Range("K74").Value = 3125
$ws.Range("L74").Value = 14848.777
$ws.Range("M74").Value = -2189
$ws.Range("N74").Value = -16720.777
$ws.Range("H77").Value = 10159.267
$ws.Range("I77").Value = 3125
$ws.Range("J77").Value = 14848.777
$ws.Range("K77").Value = 15625
$ws.Range("L77").Value = 74243.88499999999
$ws.Range("M77").Value = -10945
$ws.Range("N77").Value = -83603.88499999999
$ws.Range("H80").Value = 1485.6923
$ws.Range("J80").Value = 2179
$ws.Range("L80").Value = 6537
$ws.Range("N80").Value = -8533
$ws.Range("H83").Value = 1485.6923
$ws.Range("J83").Value = 2179
$ws.Range("L83").Value = 19611
$ws.Range("N83").Value = -29595
$ws.Range("H97").Value = 200659.6
$ws.Range("J97").Value = 200659.6
$ws.Range("L97").Value = 601978.8
$ws.Range("N97").Value = -602970.8
$ws.Range("H101").Value = 1922.5
$ws.Range("I101").Value = 1562.2
$ws.Range("J101").Value = 2222.75
$ws.Range("K101").Value = 4686.6
$ws.Range("L101").Value = 6668.25
$ws.Range("M101").Value = -3064.6
$ws.Range("N101").Value = -9912.25
$ws.Range("H106").Value = 6446.615
$ws.Range("I106").Value = 5365.3335
$ws.Range("K106").Value = 5365.3335
$ws.Range("M106").Value = -4734.3335
$ws.Range("H107").Value = 1209.15
$ws.Range("J107").Value = 2122.8333
$ws.Range("L107").Value = 2122.8333
$ws.Range("N107").Value = -5962.8333
$ws.Range("H131").Value = 1293.909
$ws.Range("I131").Value = 1268.1111
$ws.Range("J131").Value = 1410
$ws.Range("K131").Value = 3804.3333
$ws.Range("L131").Value = 4230
$ws.Range("M131").Value = 1235.6667
$ws.Range("N131").Value = -14310
$ws.Range("H132").Value = 2065.9443
$ws.Range("I132").Value = 2009.9183
$ws.Range("J132").Value = 2615
$ws.Range("K132").Value = 6029.7549
$ws.Range("L132").Value = 7845
$ws.Range("M132").Value = -3499.7549
$ws.Range("N132").Value = -12905
$ws.Range("H135").Value = 31659.725
$ws.Range("I135").Value = 33993.777
$ws.Range("K135").Value = 305943.993
$ws.Range("M135").Value = -303408.993
$ws.Range("H137").Value = 108365.695
$ws.Range("I137").Value = 2189.0435
$ws.Range("J137").Value = 395667.25
$ws.Range("K137").Value = 6567.130500000001
$ws.Range("L137").Value = 1187001.75
$ws.Range("M137").Value = -4017.130500000001
$ws.Range("N137").Value = -1192101.75
$ws.Range("H138").Value = 4139.4375
$ws.Range("I138").Value = 3752.3333
$ws.Range("J138").Value = 4290.913
$ws.Range("K138").Value = 11256.9999
$ws.Range("L138").Value = 12872.739
$ws.Range("M138").Value = -6116.999899999999
$ws.Range("N138").Value = -23152.739
$ws.Range("H141").Value = 410.9
$ws.Range("I141").Value = 402.91666
$ws.Range("J141").Value = 602.5
$ws.Range("K141").Value = 1208.74998
$ws.Range("L141").Value = 1807.5
$ws.Range("M141").Value = 3971.25002
$ws.Range("N141").Value = -12167.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1568.3334
$ws.Range("I2").Value = 1810.591
$ws.Range("J2").Value = 902.125
$ws.Range("K2").Value = 1810.591
$ws.Range("L2").Value = 902.125
$ws.Range("M2").Value = -1697.591
$ws.Range("N2").Value = -1128.125
$ws.Range("H6").Value = 2469666.2
$ws.Range("I6").Value = 2469666.2
$ws.Range("K6").Value = 2469666.2
$ws.Range("M6").Value = -2469493.2
$ws.Range("H32").Value = 12925.172
$ws.Range("I32").Value = 8876.5625
$ws.Range("K32").Value = 8876.5625
$ws.Range("M32").Value = -8589.5625
$ws.Range("H61").Value = 23582.682
$ws.Range("I61").Value = 3425.375
$ws.Range("K61").Value = 3425.375
$ws.Range("M61").Value = -3213.375
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("H63").Value = 2064.1738
$ws.Range("I63").Value = 2076.8948
$ws.Range("K63").Value = 2076.8948
$ws.Range("M63").Value = -1390.8948
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0
$ws.Range("H66").Value = 2064.1738
$ws.Range("I66").Value = 2076.8948
$ws.Range("K66").Value = 10384.474
$ws.Range("M66").Value = -6952.474
$ws.Range("H74").Value = 320712.22
$ws.Range("I74").Value = 322289.75
$ws.Range("K74").Value = 322289.75
$ws.Range("M74").Value = -321415.75
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31748
$ws.Range("H77").Value = 320712.22
$ws.Range("I77").Value = 322289.75
$ws.Range("K77").Value = 1611448.75
$ws.Range("M77").Value = -1607080.75
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -98736
$ws.Range("H82").Value = 41000
$ws.Range("J82").Value = 41000
$ws.Range("L82").Value = 41000
$ws.Range("N82").Value = -41722
$ws.Range("H85").Value = 41000
$ws.Range("J85").Value = 41000
$ws.Range("L85").Value = 41000
$ws.Range("N85").Value = -43496
$ws.Range("H102").Value = 305127.3
$ws.Range("I102").Value = 436728.2
$ws.Range("J102").Value = 2445.3
$ws.Range("K102").Value = 436728.2
$ws.Range("L102").Value = 2445.3
$ws.Range("M102").Value = -435106.2
$ws.Range("N102").Value = -5689.3
$ws.Range("H116").Value = 1568.3334
$ws.Range("I116").Value = 1810.591
$ws.Range("J116").Value = 902.125
$ws.Range("K116").Value = 1810.591
$ws.Range("L116").Value = 902.125
$ws.Range("M116").Value = 483.4090000000001
$ws.Range("N116").Value = -5490.125
$ws.Range("H122").Value = 33427.2
$ws.Range("I122").Value = 3155.476
$ws.Range("J122").Value = 104061.22
$ws.Range("K122").Value = 9466.428
$ws.Range("L122").Value = 312183.66
$ws.Range("M122").Value = -7016.428
$ws.Range("N122").Value = -317083.66
$ws.Range("H132").Value = 2314.2778
$ws.Range("I132").Value = 2222.8936
$ws.Range("J132").Value = 2927.8572
$ws.Range("K132").Value = 6668.6808
$ws.Range("L132").Value = 8783.571599999999
$ws.Range("M132").Value = -4138.6808
$ws.Range("N132").Value = -13843.5716
$ws.Range("H136").Value = 23582.682
$ws.Range("I136").Value = 3425.375
$ws.Range("K136").Value = 10276.125
$ws.Range("M136").Value = -7726.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1568.3334
$ws.Range("I3").Value = 1810.591
$ws.Range("J3").Value = 902.125
$ws.Range("K3").Value = 1810.591
$ws.Range("L3").Value = 902.125
$ws.Range("M3").Value = -1696.591
$ws.Range("N3").Value = -1130.125
$ws.Range("H5").Value = 2317.3635
$ws.Range("I5").Value = 533.3333
$ws.Range("K5").Value = 533.3333
$ws.Range("M5").Value = -420.3333
$ws.Range("H57").Value = 70780
$ws.Range("J57").Value = 70780
$ws.Range("L57").Value = 70780
$ws.Range("N57").Value = -72220
$ws.Range("H76").Value = 20750
$ws.Range("J76").Value = 20750
$ws.Range("L76").Value = 20750
$ws.Range("N76").Value = -21380
$ws.Range("H79").Value = 20750
$ws.Range("J79").Value = 20750
$ws.Range("L79").Value = 20750
$ws.Range("N79").Value = -22934
$ws.Range("H86").Value = 18216854
$ws.Range("I86").Value = 25666298
$ws.Range("K86").Value = 25666298
$ws.Range("M86").Value = -25665175
$ws.Range("H89").Value = 18216854
$ws.Range("I89").Value = 25666298
$ws.Range("K89").Value = 128331490
$ws.Range("M89").Value = -128325874
$ws.Range("H99").Value = 2177.0833
$ws.Range("I99").Value = 2185.875
$ws.Range("J99").Value = 2159.5
$ws.Range("K99").Value = 2185.875
$ws.Range("L99").Value = 2159.5
$ws.Range("M99").Value = -687.875
$ws.Range("N99").Value = -5155.5
$ws.Range("H105").Value = 23812732
$ws.Range("I105").Value = 31253154
$ws.Range("J105").Value = 3380.8
$ws.Range("K105").Value = 31253154
$ws.Range("L105").Value = 3380.8
$ws.Range("M105").Value = -31251407
$ws.Range("N105").Value = -6874.8
$ws.Range("H132").Value = 87462.5
$ws.Range("J132").Value = 87462.5
$ws.Range("L132").Value = 87462.5
$ws.Range("N132").Value = -97582.5
$ws.Range("H134").Value = 4714.6206
$ws.Range("I134").Value = 4294.619
$ws.Range("K134").Value = 12883.857
$ws.Range("M134").Value = -10348.857
$ws.Range("H136").Value = 70780
$ws.Range("J136").Value = 70780
$ws.Range("L136").Value = 70780
$ws.Range("N136").Value = -80980
$ws.Range("H137").Value = 99490
$ws.Range("J137").Value = 99490
$ws.Range("L137").Value = 99490
$ws.Range("N137").Value = -109690
$ws.Range("H138").Value = 87966.11
$ws.Range("J138").Value = 92711.875
$ws.Range("L138").Value = 92711.875
$ws.Range("N138").Value = -102991.875
$ws.Range("H140").Value = 80950.836
$ws.Range("I140").Value = 72854.5
$ws.Range("J140").Value = 84999
$ws.Range("K140").Value = 72854.5
$ws.Range("L140").Value = 84999
$ws.Range("M140").Value = -67674.5
$ws.Range("N140").Value = -95359
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3208.55
$ws.Range("I16").Value = 2781.5386
$ws.Range("K16").Value = 2781.5386
$ws.Range("M16").Value = -2494.5386
$ws.Range("H22").Value = 800.1429000000001
$ws.Range("I22").Value = 263.66666
$ws.Range("K22").Value = 263.66666
$ws.Range("M22").Value = 86.33334000000002
$ws.Range("H31").Value = 2270.0989
$ws.Range("I31").Value = 1818
$ws.Range("J31").Value = 2753.0227
$ws.Range("K31").Value = 1818
$ws.Range("L31").Value = 2753.0227
$ws.Range("M31").Value = -1523
$ws.Range("N31").Value = -3343.0227
$ws.Range("H32").Value = 6254.2856
$ws.Range("I32").Value = 2195
$ws.Range("J32").Value = 11666.667
$ws.Range("K32").Value = 2195
$ws.Range("L32").Value = 11666.667
$ws.Range("M32").Value = -1879
$ws.Range("N32").Value = -12298.667
$ws.Range("H34").Value = 2270.0989
$ws.Range("I34").Value = 1818
$ws.Range("J34").Value = 2753.0227
$ws.Range("K34").Value = 1818
$ws.Range("L34").Value = 2753.0227
$ws.Range("M34").Value = -1616
$ws.Range("N34").Value = -3157.0227
$ws.Range("H57").Value = 28749.75
$ws.Range("J57").Value = 28749.75
$ws.Range("L57").Value = 28749.75
$ws.Range("N57").Value = -29869.75
$ws.Range("H105").Value = 1576.8846
$ws.Range("I105").Value = 1222.5
$ws.Range("K105").Value = 1222.5
$ws.Range("M105").Value = 524.5
$ws.Range("H113").Value = 3208.55
$ws.Range("I113").Value = 2781.5386
$ws.Range("K113").Value = 2781.5386
$ws.Range("M113").Value = -611.5385999999999
$ws.Range("H132").Value = 5217.9697
$ws.Range("I132").Value = 2544.4443
$ws.Range("J132").Value = 17248.834
$ws.Range("K132").Value = 7633.3329
$ws.Range("L132").Value = 51746.50199999999
$ws.Range("M132").Value = -5103.3329
$ws.Range("N132").Value = -56806.50199999999
$ws.Range("H134").Value = 2235.925
$ws.Range("I134").Value = 1308.138
$ws.Range("J134").Value = 4681.909
$ws.Range("K134").Value = 3924.414
$ws.Range("L134").Value = 14045.727
$ws.Range("M134").Value = -1389.414
$ws.Range("N134").Value = -19115.727
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 692
$ws.Range("I26").Value = 600
$ws.Range("K26").Value = 1800
$ws.Range("M26").Value = -1512
$ws.Range("H52").Value = 1496
$ws.Range("J52").Value = 1496
$ws.Range("L52").Value = 4488
$ws.Range("N52").Value = -5020
$ws.Range("H80").Value = 25665.555
$ws.Range("I80").Value = 20497.834
$ws.Range("J80").Value = 36001
$ws.Range("K80").Value = 61493.50199999999
$ws.Range("L80").Value = 108003
$ws.Range("M80").Value = -60557.50199999999
$ws.Range("N80").Value = -109875
$ws.Range("H83").Value = 25665.555
$ws.Range("I83").Value = 20497.834
$ws.Range("J83").Value = 36001
$ws.Range("K83").Value = 184480.506
$ws.Range("L83").Value = 324009
$ws.Range("M83").Value = -179800.506
$ws.Range("N83").Value = -333369
$ws.Range("H86").Value = 330.8
$ws.Range("J86").Value = 326.3
$ws.Range("L86").Value = 978.9000000000001
$ws.Range("N86").Value = -3350.9
$ws.Range("H89").Value = 330.8
$ws.Range("J89").Value = 326.3
$ws.Range("L89").Value = 2936.7
$ws.Range("N89").Value = -14792.7
$ws.Range("H97").Value = 358.5
$ws.Range("J97").Value = 299.42856
$ws.Range("L97").Value = 898.28568
$ws.Range("N97").Value = -1890.28568
$ws.Range("H107").Value = 1299.88
$ws.Range("J107").Value = 1595.5264
$ws.Range("L107").Value = 4786.5792
$ws.Range("N107").Value = -8626.5792
$ws.Range("H113").Value = 1430.6086
$ws.Range("I113").Value = 676.5
$ws.Range("J113").Value = 1696.7646
$ws.Range("K113").Value = 2029.5
$ws.Range("L113").Value = 5090.293799999999
$ws.Range("M113").Value = 140.5
$ws.Range("N113").Value = -9430.293799999999
$ws.Range("H130").Value = 2806.75
$ws.Range("I130").Value = 1909
$ws.Range("K130").Value = 5727
$ws.Range("M130").Value = -707
$ws.Range("H137").Value = 3375.3125
$ws.Range("I137").Value = 1406.3846
$ws.Range("J137").Value = 11907.333
$ws.Range("K137").Value = 4219.1538
$ws.Range("L137").Value = 35721.999
$ws.Range("M137").Value = 880.8462
$ws.Range("N137").Value = -45921.999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1555.9474
$ws.Range("I3").Value = 1074.0769
$ws.Range("K3").Value = 1074.0769
$ws.Range("M3").Value = -958.0769
$ws.Range("H7").Value = 14282.286
$ws.Range("J7").Value = 14282.286
$ws.Range("L7").Value = 14282.286
$ws.Range("N7").Value = -14506.286
$ws.Range("H8").Value = 14282.286
$ws.Range("J8").Value = 14282.286
$ws.Range("L8").Value = 14282.286
$ws.Range("N8").Value = -14560.286
$ws.Range("H10").Value = 801292.6
$ws.Range("I10").Value = 1333966.6
$ws.Range("K10").Value = 1333966.6
$ws.Range("M10").Value = -1333797.6
$ws.Range("H12").Value = 1014799.4
$ws.Range("I12").Value = 1263499.8
$ws.Range("K12").Value = 1263499.8
$ws.Range("M12").Value = -1263359.8
$ws.Range("H14").Value = 4441.4707
$ws.Range("I14").Value = 944
$ws.Range("J14").Value = 8376.125
$ws.Range("K14").Value = 944
$ws.Range("L14").Value = 8376.125
$ws.Range("M14").Value = -776
$ws.Range("N14").Value = -8712.125
$ws.Range("H36").Value = 10010347
$ws.Range("I36").Value = 12738776
$ws.Range("J36").Value = 6108.3335
$ws.Range("K36").Value = 12738776
$ws.Range("L36").Value = 6108.3335
$ws.Range("M36").Value = -12738291
$ws.Range("N36").Value = -7078.3335
$ws.Range("H43").Value = 9442.24
$ws.Range("I43").Value = 4404.2666
$ws.Range("K43").Value = 4404.2666
$ws.Range("M43").Value = -4253.2666
$ws.Range("H93").Value = 85000
$ws.Range("J93").Value = 85000
$ws.Range("L93").Value = 85000
$ws.Range("N93").Value = -88744
$ws.Range("H102").Value = 27872.783
$ws.Range("I102").Value = 34562.906
$ws.Range("J102").Value = 12581.071
$ws.Range("K102").Value = 34562.906
$ws.Range("L102").Value = 12581.071
$ws.Range("M102").Value = -32940.906
$ws.Range("N102").Value = -15825.071
$ws.Range("H122").Value = 140143.08
$ws.Range("J122").Value = 2966.3333
$ws.Range("L122").Value = 8898.999899999999
$ws.Range("N122").Value = -13798.9999
$ws.Range("H126").Value = 39451.72
$ws.Range("I126").Value = 53468
$ws.Range("J126").Value = 3409.8572
$ws.Range("K126").Value = 160404
$ws.Range("L126").Value = 10229.5716
$ws.Range("M126").Value = -157934
$ws.Range("N126").Value = -15169.5716
$ws.Range("H132").Value = 5723.4116
$ws.Range("I132").Value = 5979.8667
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 17939.6001
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -15409.6001
$ws.Range("N132").Value = -16460
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 4270.2856
$ws.Range("J9").Value = 3966.4
$ws.Range("L9").Value = 3966.4
$ws.Range("N9").Value = -4414.4
$ws.Range("H16").Value = 1608.3125
$ws.Range("I16").Value = 1608.3125
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1608.3125
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1438.3125
$ws.Range("H20").Value = 16749.75
$ws.Range("J20").Value = 7000
$ws.Range("L20").Value = 7000
$ws.Range("N20").Value = -7452
$ws.Range("H22").Value = 3125.66
$ws.Range("J22").Value = 1925.3939
$ws.Range("L22").Value = 1925.3939
$ws.Range("N22").Value = -2515.3939
$ws.Range("H27").Value = 3125.66
$ws.Range("J27").Value = 1925.3939
$ws.Range("L27").Value = 1925.3939
$ws.Range("N27").Value = -2139.3939
$ws.Range("H40").Value = 6003
$ws.Range("I40").Value = 5208.143
$ws.Range("K40").Value = 5208.143
$ws.Range("M40").Value = -5072.143
$ws.Range("H43").Value = 33466.25
$ws.Range("I43").Value = 28138.584
$ws.Range("K43").Value = 28138.584
$ws.Range("M43").Value = -27945.584
$ws.Range("H46").Value = 3046.4285
$ws.Range("I46").Value = 583
$ws.Range("K46").Value = 583
$ws.Range("M46").Value = -395
$ws.Range("H61").Value = 47620884
$ws.Range("I61").Value = 50001680
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 50001680
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -50001478
$ws.Range("N61").Value = -5404
$ws.Range("H62").Value = 847631.4
$ws.Range("J62").Value = 847631.4
$ws.Range("L62").Value = 847631.4
$ws.Range("N62").Value = -848879.4
$ws.Range("H65").Value = 847631.4
$ws.Range("J65").Value = 847631.4
$ws.Range("L65").Value = 2542894.2
$ws.Range("N65").Value = -2549134.2
$ws.Range("H76").Value = 15392
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 15392
$ws.Range("K76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("M76").Value = 15392
$ws.Range("N76").Value = -16068
$ws.Range("H79").Value = 15392
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 15392
$ws.Range("K79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("M79").Value = 15392
$ws.Range("N79").Value = -17732
$ws.Range("H113").Value = 47620884
$ws.Range("I113").Value = 50001680
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 50001680
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -49999510
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 3617.7632
$ws.Range("I122").Value = 3541
$ws.Range("K122").Value = 10623
$ws.Range("M122").Value = -8173
$ws.Range("H132").Value = 3665.8955
$ws.Range("I132").Value = 2550.673
$ws.Range("J132").Value = 7532
$ws.Range("K132").Value = 7652.018999999999
$ws.Range("L132").Value = 22596
$ws.Range("M132").Value = -5122.018999999999
$ws.Range("N132").Value = -27656
$ws.Range("H133").Value = 91402
$ws.Range("I133").Value = 49596
$ws.Range("J133").Value = 99763.2
$ws.Range("K133").Value = 49596
$ws.Range("L133").Value = 99763.2
$ws.Range("M133").Value = -47066
$ws.Range("N133").Value = -104823.2
$ws.Range("H136").Value = 41286.293
$ws.Range("I136").Value = 3355.3125
$ws.Range("K136").Value = 10065.9375
$ws.Range("M136").Value = -7515.9375
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = 0
$ws.Range("H14").Value = 4310.143
$ws.Range("I14").Value = 3318.2856
$ws.Range("J14").Value = 5302
$ws.Range("K14").Value = 3318.2856
$ws.Range("L14").Value = 5302
$ws.Range("M14").Value = -3150.2856
$ws.Range("N14").Value = -5638
$ws.Range("H81").Value = 1624.25
$ws.Range("I81").Value = 1713.4286
$ws.Range("K81").Value = 3426.8572
$ws.Range("M81").Value = -2365.8572
$ws.Range("H84").Value = 1624.25
$ws.Range("I84").Value = 1713.4286
$ws.Range("K84").Value = 17134.286
$ws.Range("M84").Value = -11830.286
$ws.Range("H132").Value = 1581.75
$ws.Range("I132").Value = 1522.1428
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 4566.428400000001
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -2036.428400000001
$ws.Range("N132").Value = -11057
$ws.Range("H136").Value = 6488.3945
$ws.Range("I136").Value = 7774.7085
$ws.Range("J136").Value = 4283.2856
$ws.Range("K136").Value = 23324.1255
$ws.Range("L136").Value = 12849.8568
$ws.Range("M136").Value = -20774.1255
$ws.Range("N136").Value = -17949.8568
